# Milestone 5 file dump
# Adds two new sound-effect rows (shriek on death / victory bells), fills in
# the missing "Category" (B) column for the wood/metal collision rows, and
# repoints the "metal collision" asset link at the Gammelsmurfen778 sound
# (replacing the old launemax link) with a live hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20 (wood collision): fill in the Category column ---
$ws.Range("B20").Value = "sound effect"
$ws.Range("B2").Copy()
$ws.Range("B20").PasteSpecial(-4122)

# --- Row 21 (metal collision): fill in Category + replace the asset link ---
$ws.Range("B21").Value = "sound effect"
$ws.Range("B2").Copy()
$ws.Range("B21").PasteSpecial(-4122)

$ws.Range("D21").Value = "https://freesound.org/people/Gammelsmurfen778/sounds/474007/"
$ws.Hyperlinks.Add($ws.Range("D21"), "https://freesound.org/people/Gammelsmurfen778/sounds/474007/")
$ws.Range("D6").Copy()
$ws.Range("D21").PasteSpecial(-4122)

# --- Row 22 (new): shriek (upon death) ---
$ws.Range("D22").Value = "https://freesound.org/people/JohnsonBrandEditing/sounds/173944/"

$ws.Range("A22").Value = "shriek (upon death)"
$ws.Range("A20").Copy()
$ws.Range("A22").PasteSpecial(-4122)

$ws.Range("B22").Value = "sound effect"
$ws.Range("B2").Copy()
$ws.Range("B22").PasteSpecial(-4122)

$ws.Rows.Item(22).RowHeight = 15.75

# --- Row 23 (new): Victory bells ---
$ws.Range("A23").Value = "Victory bells"
$ws.Range("A20").Copy()
$ws.Range("A23").PasteSpecial(-4122)

$ws.Range("B23").Value = "sound effect"
$ws.Range("B2").Copy()
$ws.Range("B23").PasteSpecial(-4122)

$ws.Range("D23").Value = "https://freesound.org/people/loopsamples.club/sounds/483383/"
$ws.Range("C23").Value = "The player collects bells to win, so the victory chime is the sound of several bells."

$ws.Rows.Item(23).RowHeight = 15.75

# Clear the clipboard marquee and leave the selection where the author left it.
$excel.CutCopyMode = $false
$ws.Range("C23").Select()
